# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Vega Monumental Concepción" / Frutilla
# right above the existing row 567 (most recent date first), pushing the
# previously-existing rows 567-579 down to 569-581.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 567 (each Insert() pushes row 567+ down by one,
# and inherits formatting, e.g. number format, from the row above - same as
# Excel's native "Insert Copied/Row Above" behaviour).
$ws.Rows.Item(567).Insert()
$ws.Rows.Item(567).Insert()

# New row 567: Primera, Provincia de Melipilla
$ws.Cells.Item(567, 1).Value2 = 11
$ws.Cells.Item(567, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(567, 3).Value2 = "Bíobío"
$ws.Cells.Item(567, 4).Value2 = 45239
$ws.Cells.Item(567, 5).Value2 = 8
$ws.Cells.Item(567, 6).Value2 = "Fruta"
$ws.Cells.Item(567, 7).Value2 = 100101
$ws.Cells.Item(567, 8).Value2 = "Berries"
$ws.Cells.Item(567, 9).Value2 = 100112025
$ws.Cells.Item(567, 10).Value2 = "Frutilla"
$ws.Cells.Item(567, 11).Value2 = "Sin especificar"
$ws.Cells.Item(567, 12).Value2 = "Primera"
$ws.Cells.Item(567, 13).Value2 = 300
$ws.Cells.Item(567, 14).Value2 = 12000
$ws.Cells.Item(567, 15).Value2 = 13000
$ws.Cells.Item(567, 16).Value2 = 12500
$ws.Cells.Item(567, 17).Value2 = "`$/bandeja 7 kilos"
$ws.Cells.Item(567, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(567, 19).Value2 = 1786
$ws.Cells.Item(567, 20).Value2 = 7

# New row 568: Segunda, Provincia de Melipilla
$ws.Cells.Item(568, 1).Value2 = 11
$ws.Cells.Item(568, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(568, 3).Value2 = "Bíobío"
$ws.Cells.Item(568, 4).Value2 = 45239
$ws.Cells.Item(568, 5).Value2 = 8
$ws.Cells.Item(568, 6).Value2 = "Fruta"
$ws.Cells.Item(568, 7).Value2 = 100101
$ws.Cells.Item(568, 8).Value2 = "Berries"
$ws.Cells.Item(568, 9).Value2 = 100112025
$ws.Cells.Item(568, 10).Value2 = "Frutilla"
$ws.Cells.Item(568, 11).Value2 = "Sin especificar"
$ws.Cells.Item(568, 12).Value2 = "Segunda"
$ws.Cells.Item(568, 13).Value2 = 200
$ws.Cells.Item(568, 14).Value2 = 10000
$ws.Cells.Item(568, 15).Value2 = 10000
$ws.Cells.Item(568, 16).Value2 = 10000
$ws.Cells.Item(568, 17).Value2 = "`$/bandeja 7 kilos"
$ws.Cells.Item(568, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(568, 19).Value2 = 1429
$ws.Cells.Item(568, 20).Value2 = 7
